$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Table 1 (Problem 1, rows 5-15): fill in "Path Length OK?" column J with
# Yes/No verdicts.
# ---------------------------------------------------------------------------
$ws.Range("J6").Value = "Yes"
$ws.Range("J7").Value = "No"
$ws.Range("J8").Value = "Yes"
$ws.Range("J9").Value = "Yes"
$ws.Range("J10").Value = "Yes"
$ws.Range("J11").Value = "No"
$ws.Range("J12").Value = "Yes"
$ws.Range("J13").Value = "Yes"
$ws.Range("J14").Value = "Yes"
$ws.Range("J15").Value = "Yes"

# Highlight the "Greedy Best First Graph search with h_1" row green.
$ws.Range("D12:J12").Interior.Color = 5296274

# ---------------------------------------------------------------------------
# Table 2 (Problem 2, rows 21-31): fill in the whole results table.
# ---------------------------------------------------------------------------
$ws.Range("E22").Value = 40.392
$ws.Range("F22").Value = 9
$ws.Range("G22").Value = 3343
$ws.Range("H22").Value = 4609
$ws.Range("I22").Value = 30509
$ws.Range("J22").Value = "No"

$ws.Range("E24").Value = 4.257
$ws.Range("F24").Value = 619
$ws.Range("G24").Value = 624
$ws.Range("H24").Value = 625
$ws.Range("I24").Value = 5602
$ws.Range("J24").Value = "No"

$ws.Range("E26").Value = 14.981
$ws.Range("F26").Value = 9
$ws.Range("G26").Value = 4849
$ws.Range("H26").Value = 4851
$ws.Range("I26").Value = 44001
$ws.Range("J26").Value = "No"

$ws.Range("E28").Value = 2.358
$ws.Range("F28").Value = 16
$ws.Range("G28").Value = 966
$ws.Range("H28").Value = 968
$ws.Range("I28").Value = 8694
$ws.Range("J28").Value = "Yes"

$ws.Range("E29").Value = 11.917
$ws.Range("F29").Value = 9
$ws.Range("G29").Value = 4849
$ws.Range("H29").Value = 4851
$ws.Range("I29").Value = 44001
$ws.Range("J29").Value = "No"

$ws.Range("E30").Value = 4.239
$ws.Range("F30").Value = 9
$ws.Range("G30").Value = 1443
$ws.Range("H30").Value = 1445
$ws.Range("I30").Value = 13234
$ws.Range("J30").Value = "No"

$ws.Range("E31").Value = 48.441
$ws.Range("F31").Value = 9
$ws.Range("G31").Value = 85
$ws.Range("H31").Value = 87
$ws.Range("I31").Value = 831
$ws.Range("J31").Value = "No"

# Rows that could not be run (stack overflow / out of memory) are marked "NA"
# and highlighted orange.
$ws.Range("D23:J23").Value = "NA"
$ws.Range("D23:J23").Interior.Color = 6126847

$ws.Range("D25:J25").Value = "NA"
$ws.Range("D25:J25").Interior.Color = 6126847

$ws.Range("D27:J27").Value = "NA"
$ws.Range("D27:J27").Interior.Color = 6126847

# Highlight the "Greedy Best First Graph search with h_1" row green.
$ws.Range("D28:J28").Interior.Color = 5296274

# ---------------------------------------------------------------------------
# Restore the view state recorded in the saved workbook.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Range("G17").Select()
